$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("general")
$ws.Range("B3").Value = 147.5968802676175
$ws.Range("B4").Value = 0.01699995994567871
$ws.Range("B6").Value = 33.21688026761753
$ws.Range("B9").Value = 0
$ws.Range("B10").Value = 114.38

$ws = $wb.Worksheets.Item("x")
$ws.Range("B2").Value = 3
$ws.Range("B3").Value = 7
$ws.Range("B4").Value = 9
$ws.Range("B6").Value = 11
$ws.Range("B7").Value = 6
$ws.Range("B8").Value = 12
$ws.Range("B10").Value = 5
$ws.Range("B12").Value = 4
$ws.Range("B13").Value = 10
$ws.Range("B14").Value = 8

$ws = $wb.Worksheets.Item("U")
$ws.Range("B10").Value = 2
$ws.Range("B11").Value = 3

$ws = $wb.Worksheets.Item("TBar")
$ws.Range("B3").Value = 20.60033324079215
$ws.Range("B4").Value = 24.72107346555759
$ws.Range("B5").Value = 20
$ws.Range("B6").Value = 25.22705701268762
$ws.Range("B7").Value = 24.16886835983306
$ws.Range("B8").Value = 20.34885527085025
$ws.Range("B9").Value = 20
$ws.Range("B10").Value = 22.91480702620259
$ws.Range("B11").Value = 10
$ws.Range("B13").Value = 27.3015579161986
$ws.Range("B14").Value = 30
$ws.Range("B15").Value = 25.35398438790795

$ws = $wb.Worksheets.Item("Q")
$ws.Range("C7").Value = 193.0200000000017
$ws.Range("C8").Value = 202.3100000000017
$ws.Range("C9").Value = 191.2450000000017
$ws.Range("C10").Value = 208.9250000000017
$ws.Range("C11").Value = 197.6600000000017
$ws.Range("C12").Value = 274.4950000000024
$ws.Range("C13").Value = 282.9900000000024
$ws.Range("C14").Value = 275.9600000000024
$ws.Range("C15").Value = 289.3600000000025
$ws.Range("C16").Value = 285.0050000000024
$ws.Range("C17").Value = 46.91999999999942
$ws.Range("C18").Value = 36.10499999999942
$ws.Range("C19").Value = 34.91499999999942
$ws.Range("C20").Value = 37.48999999999942
$ws.Range("C21").Value = 39.43499999999941
$ws.Range("C22").Value = 272.2599999999987
$ws.Range("C23").Value = 291.1899999999987
$ws.Range("C24").Value = 278.0049999999987
$ws.Range("C25").Value = 287.35
$ws.Range("C26").Value = 263.9399999999987
$ws.Range("C27").Value = 224.1799999999995
$ws.Range("C28").Value = 224.6649999999995
$ws.Range("C29").Value = 201.1149999999995
$ws.Range("C30").Value = 218.9699999999995
$ws.Range("C31").Value = 207.1049999999995
$ws.Range("C32").Value = 154.3
$ws.Range("C33").Value = 148.3449999999993
$ws.Range("C34").Value = 128.7049999999993
$ws.Range("C35").Value = 146.3249999999992
$ws.Range("C36").Value = 134.2149999999993
$ws.Range("C37").Value = 83.07500000000087
$ws.Range("C38").Value = 84.72000000000087
$ws.Range("C39").Value = 77.97500000000086
$ws.Range("C40").Value = 89.77000000000088
$ws.Range("C41").Value = 81.85500000000087
$ws.Range("C42").Value = 199.6299999999992
$ws.Range("C43").Value = 211.1349999999992
$ws.Range("C44").Value = 195.6849999999992
$ws.Range("C45").Value = 203.3949999999992
$ws.Range("C46").Value = 183.1249999999992
$ws.Range("C47").Value = 73.75500000000051
$ws.Range("C48").Value = 73.2950000000005
$ws.Range("C49").Value = 74.1350000000005
$ws.Range("C50").Value = 78.97000000000051
$ws.Range("C51").Value = 71.28000000000051
$ws.Range("C52").Value = 131.7049999999998
$ws.Range("C53").Value = 131.9749999999998
$ws.Range("C54").Value = 135.8649999999998
$ws.Range("C55").Value = 139.6249999999998
$ws.Range("C56").Value = 124.1949999999998
$ws.Range("C57").Value = 274.4950000000024
$ws.Range("C58").Value = 282.9900000000024
$ws.Range("C59").Value = 275.9600000000024
$ws.Range("C60").Value = 289.3600000000025
$ws.Range("C61").Value = 285.0050000000024
$ws.Range("C62").Value = 272.2599999999987
$ws.Range("C63").Value = 291.1899999999987
$ws.Range("C64").Value = 278.0049999999987
$ws.Range("C65").Value = 287.35
$ws.Range("C66").Value = 263.9399999999987
$ws.Range("C67").Value = 224.1799999999995
$ws.Range("C68").Value = 224.6649999999995
$ws.Range("C69").Value = 201.1149999999995
$ws.Range("C70").Value = 218.9699999999995
$ws.Range("C71").Value = 207.1049999999995

$ws = $wb.Worksheets.Item("R")
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("C15").Value = 0

$ws = $wb.Worksheets.Item("L")
$ws.Range("C42").Value = 11.73
$ws.Range("C43").Value = 14.67
$ws.Range("C44").Value = 5.58
$ws.Range("C45").Value = 11.76
$ws.Range("C46").Value = 13.45
$ws.Range("C47").Value = 0
$ws.Range("C48").Value = 0
$ws.Range("C49").Value = 0
$ws.Range("C50").Value = 0
$ws.Range("C51").Value = 0

# rho sheet: clear all data rows except header, per diff (dimension becomes A1:C1)
$ws = $wb.Worksheets.Item("rho")
$ws.Rows("2:8").Delete()
